$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $escaped = $text -replace '"', '""'
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

$ws.Range('D2').Value = '42.399.07'
$ws.Range('E2').Value = '  +1.52%  '
$ws.Range('D3').Value = '2.279.01'
$ws.Range('E3').Value = '  +0.50%  '
Set-TextValue $ws.Range('D5') '307.29'
$ws.Range('E5').Value = '  +1.33%  '
Set-TextValue $ws.Range('D6') '97.55'
$ws.Range('E6').Value = '  +5.83%  '
$ws.Range('E7').Value = '  +0.13%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('E9').Value = '  +2.24%  '
Set-TextValue $ws.Range('D10') '35.81'
$ws.Range('E10').Value = '  +10.81%  '
Set-TextValue $ws.Range('D11') '0.0797'
$ws.Range('E11').Value = '  +0.12%  '
$ws.Range('E12').Value = '  -0.98%  '
$ws.Range('E13').Value = '  +0.61%  '
$ws.Range('D14').Value = '2.632.67'
$ws.Range('E14').Value = '  +0.54%  '
$ws.Range('E15').Value = '  +1.84%  '
$ws.Range('D16').Value = '2.280.93'
$ws.Range('E16').Value = '  +0.61%  '
Set-TextValue $ws.Range('D17') '0.802'
$ws.Range('E17').Value = '  +3.98%  '
$ws.Range('D18').Value = '42.315.67'
$ws.Range('E18').Value = '  +1.56%  '
Set-TextValue $ws.Range('D19') '12.60'
$ws.Range('E19').Value = '  +0.96%  '
$ws.Range('D20').Value = '0.0₃0911'
$ws.Range('E20').Value = '  +0.79%  '
$ws.Range('E21').Value = '  +0.81%  '
Set-TextValue $ws.Range('D22') '67.71'
$ws.Range('E22').Value = '  +1.06%  '
Set-TextValue $ws.Range('D23') '241.06'
$ws.Range('E23').Value = '  +0.60%  '
$ws.Range('E24').Value = '  +0.10%  '
$ws.Range('E25').Value = '  +1.85%  '
$ws.Range('E26').Value = '  -0.02%  '
Set-TextValue $ws.Range('D27') '23.90'
$ws.Range('E27').Value = '  -0.21%  '
Set-TextValue $ws.Range('D28') '37.65'
$ws.Range('E28').Value = '  +6.05%  '
Set-TextValue $ws.Range('D29') '9.52'
$ws.Range('E29').Value = '  -0.05%  '
$ws.Range('E30').Value = '  +1.04%  '
Set-TextValue $ws.Range('D31') '159.77'
$ws.Range('E31').Value = '  -0.34%  '
$ws.Range('E32').Value = '  +0.49%  '
Set-TextValue $ws.Range('D33') '0.999'
Set-TextValue $ws.Range('D34') '3.14'
$ws.Range('E35').Value = '  -0.29%  '
Set-TextValue $ws.Range('D36') '17.01'
$ws.Range('E36').Value = '  +0.59%  '
Set-TextValue $ws.Range('D37') '2.38'
$ws.Range('E37').Value = '  +0.54%  '
$ws.Range('E38').Value = '  +1.24%  '
Set-TextValue $ws.Range('D39') '1.84'
$ws.Range('E39').Value = '  +2.93%  '
$ws.Range('E40').Value = '  -0.96%  '
Set-TextValue $ws.Range('D41') '4.12'
$ws.Range('E41').Value = '  +5.46%  '
$ws.Range('E42').Value = '  +14.39%  '
$ws.Range('D43').Value = '1.998.67'
$ws.Range('E43').Value = '  -0.20%  '
$ws.Range('E44').Value = '  +2.17%  '
Set-TextValue $ws.Range('D45') '19.00'
$ws.Range('E45').Value = '  -3.01%  '
Set-TextValue $ws.Range('D46') '2.98'
$ws.Range('E46').Value = '  +3.48%  '
Set-TextValue $ws.Range('D47') '10.01'
$ws.Range('E47').Value = '  -2.94%  '
Set-TextValue $ws.Range('D48') '52.98'
$ws.Range('E48').Value = '  +1.21%  '
$ws.Range('E49').Value = '  +0.69%  '
Set-TextValue $ws.Range('D50') '72.07'
$ws.Range('E50').Value = '  +0.20%  '
Set-TextValue $ws.Range('D51') '92.05'
$ws.Range('E51').Value = '  +1.16%  '

$excel.CutCopyMode = 0
